$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe (Excel keeps them as text)
$ws.Range("D2").Value = "42.803.33"
$ws.Range("D3").Value = "2.260.60"
$ws.Range("E3").Value = "  +4.16%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +7.17%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +16.98%  "
$ws.Range("E10").Value = "  +9.65%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  +4.83%  "
$ws.Range("E13").Value = "  +9.06%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("D15").Value = "2.590.31"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("E16").Value = "  +4.55%  "
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "2.259.64"
$ws.Range("E18").Value = "  +4.05%  "
$ws.Range("D19").Value = "42.741.76"
$ws.Range("E19").Value = "  +4.61%  "
$ws.Range("D20").Value = "0.0₃0993"
$ws.Range("E20").Value = "  +6.45%  "
$ws.Range("E21").Value = "  +3.63%  "
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("E24").Value = "  -0.55%  "
$ws.Range("E25").Value = "  +6.79%  "
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("E33").Value = "  +15.44%  "
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("E35").Value = "  +6.55%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("E36").Value = "  +26.75%  "
$ws.Range("E37").Value = "  +4.07%  "
$ws.Range("E38").Value = "  +12.29%  "
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("E40").Value = "  +5.60%  "
$ws.Range("E41").Value = "  +6.52%  "
$ws.Range("E42").Value = "  +9.91%  "
$ws.Range("E43").Value = "  +6.67%  "
$ws.Range("E44").Value = "  +3.25%  "
$ws.Range("E45").Value = "  +8.15%  "
$ws.Range("E46").Value = "  +5.50%  "
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("E48").Value = "  +3.55%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("E50").Value = "  +3.84%  "
$ws.Range("E51").Value = "  +4.58%  "

# Values that look like numbers (e.g. "249.49") must be forced to remain text,
# matching the original inlineStr/text cell type in the workbook.
$numericLookingValues = @{
    "D5" = "249.49"
    "D6" = "0.626"
    "D7" = "70.92"
    "D9" = "0.660"
    "D10" = "38.98"
    "D11" = "59.67"
    "D12" = "0.0971"
    "D14" = "0.104"
    "D16" = "14.90"
    "D17" = "0.880"
    "D22" = "73.10"
    "D23" = "233.72"
    "D24" = "2.08"
    "D25" = "3.94"
    "D26" = "11.48"
    "D31" = "167.88"
    "D32" = "20.93"
    "D35" = "0.0796"
    "D36" = "31.22"
    "D37" = "0.126"
    "D38" = "4.44"
    "D39" = "4.73"
    "D42" = "12.49"
    "D43" = "5.82"
    "D44" = "62.26"
    "D45" = "9.10"
    "D47" = "4.83"
    "D49" = "1.01"
}
foreach ($cellRef in $numericLookingValues.Keys) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $numericLookingValues[$cellRef]
    $rng.Style = "Normal"
}
